$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Conditional formatting on C78:D78 doesn't need to recalc on every edit.
$ws.EnableFormatConditionsCalculation = $false

# --- Section 1 test table (rows 18-34): enter measured resistance values ---
$ws.Range("K18").Value = "n/a"
$ws.Range("K19").Value = 0.7
$ws.Range("K20").Value = 0.7
$ws.Range("K21").Value = 0.7
$ws.Range("K22").Value = 0.7
$ws.Range("K23").Value = 0.7
$ws.Range("K24").Value = 0.7
$ws.Range("K25").Value = 0.7
$ws.Range("K26").Value = 0.8
$ws.Range("K27").Value = 0.7
$ws.Range("K28").Value = 0.7
$ws.Range("K29").Value = 0.7
$ws.Range("K30").Value = 0.7
$ws.Range("K31").Value = 0.7
$ws.Range("K32").Value = 0.8
$ws.Range("K33").Value = 0.7

# Row 34 previously had no ok/NOK or fail-flag formulas next to its limit value (K34).
# Extend the same pattern used for K18:K33 down into row 34.
$ws.Range("L34").Formula = '=IF(AND(K34<K$34,ISNUMBER(K34)),"ok","NOK")'
$ws.Range("P34").Formula = '=IF(L34="NOK",1,0)'

# --- Section 2 test table (rows 40-47): enter measured values ---
$ws.Range("K40").Value = 0.3
$ws.Range("K41").Value = 0.4
$ws.Range("K42").Value = 0
$ws.Range("K43").Value = 0.3
$ws.Range("K44").Value = 0.3
$ws.Range("K45").Value = 0
$ws.Range("K46").Value = 0.3
$ws.Range("K47").Value = 0.3

# --- "4. Other tests" visual-inspection checks: n -> y ---
$ws.Range("M57").Value = "y"
$ws.Range("M58").Value = "y"
$ws.Range("M59").Value = "y"

# --- LED test result ---
$ws.Range("B60").Value = "pass"

# --- HV test: voltage/current readings ---
$ws.Range("K66").Value = 98.8
$ws.Range("K67").Value = 9.878

# --- Voltage drop measurements ---
$ws.Range("C70").Value = -0.0169
$ws.Range("C71").Value = -0.0165

# --- Final decision: tester name and test date ---
$ws.Range("B80").Value = "Amanda"
$ws.Range("F80").Value = "10/6/2014"

# Restore the selection to reflect where the user finished data entry.
$ws.Range("F81").Select()
